# GSC export update: append two more days (2025-11-08 and 2025-11-09) of
# data to the "Chart" sheet (sheet1). The "Critical issues" / "Non-critical
# issues" sheets keep their existing header row ("Issue" / "Validation" /
# "Items") - no changes needed there since the underlying text is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Helper: write a date-looking string into column A as plain text (not an
# Excel date serial). We build it via a formula that evaluates to the
# literal string, then collapse the formula down to its static value with
# a values-only paste. This avoids Excel's "looks like a date" auto
# conversion (and the accompanying number-format/style churn) that a
# direct .Value assignment of a date-shaped string would trigger.
function Set-TextValue($cell, [string]$text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# Row 35: 2025-11-08
Set-TextValue $ws.Range("A35") "2025-11-08"
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 82

# Row 36: 2025-11-09
Set-TextValue $ws.Range("A36") "2025-11-09"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 76
